$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CONTINUE_MAIN_TEST (row 23): the "start the experiment" wording is now shown
# before PART2 as well (it used to be shown only once, before PART1), so the
# text is reworded from an instruction to begin "the experiment" a single
# time to something that can be reused.
$ws.Range("B23").Value = "Jetzt beginnt das Experiment. <br> Viel Vergnügen!"
$ws.Range("C23").Value = "Now the experiment starts. <br> Have fun!"

# PART1_STIMULUS_DESCRIPTION (row 36): stop hard-coding "einige/some" (a few)
# and instead parametrize the count of stimuli with a template placeholder.
$ws.Range("B36").Value = "{{num_items}} kurze Klavierstücke"
$ws.Range("C36").Value = "{{num_items}} piano pieces"

# Update the remembered cursor/scroll position for the sheet, matching where
# the edits were made (down near PART2_STIMULUS_DESCRIPTION).
$ws.Range("C37").Select()
